$d = $word.ActiveDocument

# Locate the "abwmaps.com" text in the contact-info table (the right-aligned
# paragraph under the phone number) and turn it into a live hyperlink to the
# author's personal site.
$target = $d.Content
$target.Find.Execute("abwmaps.com")
$link = $d.Hyperlinks.Add($target, "https://abwmaps.com/")

# Re-find the now-hyperlinked text so it can be split into three runs
# ("abwmap" | "s" | ".com"), mirroring the original author's edit, while every
# run keeps the Hyperlink character style.
$full = $d.Content
$full.Find.Execute("abwmaps.com")
$start = $full.Start

# Toggling a character property on and back off forces Word to break the run
# at that boundary without changing the saved formatting.
$part1 = $d.Range($start, $start + 6)
$part1.Font.Bold = 1
$part1.Font.Bold = 0

$part2 = $d.Range($start + 6, $start + 7)
$part2.Font.Bold = 1
$part2.Font.Bold = 0

Write-Host "Linked abwmaps.com to the resume."
